$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "19-Jan" column (V) ------------------------------------------
# Copy formatting from the previous day's column (U) so the new cells reuse
# the existing style records (date format for the header, list-entry format
# for the data rows) instead of creating new ones, then overwrite the value.

$ws.Range("U1").Copy($ws.Range("V1"))
$ws.Range("V1").Value = 45310      # 19-Jan-2024

$ws.Range("U2").Copy($ws.Range("V2"))
$ws.Range("V2").Value = "Present"

$ws.Range("U3").Copy($ws.Range("V3"))
$ws.Range("V3").Value = "Absent"

$ws.Range("U4").Copy($ws.Range("V4"))
$ws.Range("V4").Value = "Absent"

# Row 5 (Dipti Shelavane) does not get a 19-Jan entry yet - the attendance
# tracked for that row currently stops at column O, so the old trailing
# P5:U5 entries are removed entirely (not merely cleared).
$ws.Range("P5:U5").Clear()

# --- Data validation dropdown ----------------------------------------------
# Previously one rectangular block C2:U5. Now the valid block is C2:V4 plus
# the shorter C5:P5 on the last row. Rebuild it as a single multi-area
# validation covering exactly C2:V4 and C5:P5.
$ws.Range("C2:U5").Validation.Delete()
$ws.Range("C2:V5").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')
$ws.Range("Q5:V5").Validation.Delete()

# --- View state -------------------------------------------------------------
# Best-effort restore of the scrolled/selected state captured in the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("O15,P8").Select()
